# Aggiornamento dati LamaMocogno al 23 agosto 2021
# Adds rows 344..357 (dates 2021-08-10 .. 2021-08-23, serials 44418..44431)
# following on from the existing data that ends at row 343 (serial 44417).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows' data: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$newData = @(
    @(44418, 0, 0, 0),
    @(44419, 0, 0, 0),
    @(44420, 0, 0, 0),
    @(44421, 0, 0, 0),
    @(44422, 0, 0, 0),
    @(44423, 0, 0, 0),
    @(44424, 1, 1, 37.46721618583739),
    @(44425, 3, 4, 149.8688647433496),
    @(44426, 0, 4, 149.8688647433496),
    @(44427, 0, 4, 149.8688647433496),
    @(44428, 1, 5, 187.3360809291869),
    @(44429, 0, 5, 187.3360809291869),
    @(44430, 1, 6, 224.8032971150243),
    @(44431, 1, 6, 224.8032971150243)
)

$startRow = 344

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i

    # Copy the formatting (style, number format, font, borders, etc.) from the
    # last populated row (343) down onto the new row before writing values,
    # so the new cells inherit the same look (date format on col A, etc.).
    $ws.Range("A343:D343").Copy($ws.Range("A$r`:D$r"))

    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = $false
